$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chapter 8 index-term coverage (column I) for chapters 6/7 rows + new R-functions rows
$ws.Range("I3").Value = "check"
$ws.Range("I4").Value = "check"
$ws.Range("I5").Value = "check"
$ws.Range("I6").Value = "check"
$ws.Range("I7").Value = "n/a"
$ws.Range("I8").Value = "n/a"
$ws.Range("I9").Value = "check"
$ws.Range("I10").Value = "n/a"
$ws.Range("I11").Value = "n/a"
$ws.Range("I12").Value = "check"
$ws.Range("I13").Value = "n/a"
$ws.Range("I14").Value = "check"
$ws.Range("I15").Value = "n/a"
$ws.Range("I16").Value = "check"
$ws.Range("I17").Value = "n/a"
$ws.Range("I18").Value = "check"
$ws.Range("I19").Value = "n/a"
$ws.Range("I20").Value = "n/a"
$ws.Range("I21").Value = "n/a"
$ws.Range("I22").Value = "n/a"
$ws.Range("I23").Value = "n/a"
$ws.Range("I24").Value = "n/a"
$ws.Range("I25").Value = "n/a"
$ws.Range("I26").Value = "check"
$ws.Range("I27").Value = "n/a"
$ws.Range("I28").Value = "n/a"
$ws.Range("I29").Value = "check"
$ws.Range("I30").Value = "check"
$ws.Range("I31").Value = "n/a"
$ws.Range("I32").Value = "n/a"
$ws.Range("I33").Value = "n/a"
$ws.Range("I34").Value = "check"
$ws.Range("I35").Value = "check"
$ws.Range("I36").Value = "check"
$ws.Range("I37").Value = "check"
$ws.Range("I38").Value = "n/a"
$ws.Range("I39").Value = "check"
$ws.Range("I40").Value = "check"
$ws.Range("I41").Value = "check"
$ws.Range("I42").Value = "check"
$ws.Range("I43").Value = "check"
$ws.Range("I44").Value = "check"
$ws.Range("I45").Value = "check"
$ws.Range("I46").Value = "n/a"
$ws.Range("I47").Value = "n/a"
$ws.Range("I48").Value = "n/a"
$ws.Range("I49").Value = "n/a"
$ws.Range("I50").Value = "n/a"
$ws.Range("I51").Value = "n/a"
$ws.Range("I52").Value = "n/a"
$ws.Range("I53").Value = "check"
$ws.Range("I54").Value = "n/a"
$ws.Range("I55").Value = "check"
$ws.Range("I56").Value = "n/a"
$ws.Range("I57").Value = "check"
$ws.Range("I58").Value = "check"
$ws.Range("I59").Value = "n/a"
$ws.Range("I60").Value = "check"
$ws.Range("I61").Value = "check"
$ws.Range("I62").Value = "n/a"
$ws.Range("I63").Value = "n/a"
$ws.Range("I64").Value = "n/a"
$ws.Range("I67").Value = " n/a"
$ws.Range("I68").Value = "check"
$ws.Range("I69").Value = "n/a"
$ws.Range("I70").Value = "check"
$ws.Range("I71").Value = "check"
$ws.Range("I72").Value = "check"
$ws.Range("I73").Value = "n/a"
$ws.Range("I74").Value = "n/a"
$ws.Range("I75").Value = "n/a"
$ws.Range("I76").Value = "check"
$ws.Range("I77").Value = "n/a"
$ws.Range("I78").Value = "n/a"
$ws.Range("I79").Value = "n/a"
$ws.Range("I80").Value = "n/a"

# Restore the final selection the author left the sheet in
$ws.Range("I81").Select()
